$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row (20134588 / JHNSON BB CRM M/RC50) was removed from the
# sheet; delete the entire row and let the remaining rows shift up.
$ws.Rows(95).Delete()
